# Actualización desde MV -datos-
# Append new daily UF bond-rate rows (16-09-2021 .. 29-09-2021) to sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 179; Date = "16-09-2021"; B = $null;  C = -0.75; D = 1.17; E = 1.91; F = $null;  G = 2.5 },
    @{ Row = 180; Date = "20-09-2021"; B = $null;  C = -0.74; D = 1.17; E = 1.91; F = 2.44;  G = $null },
    @{ Row = 181; Date = "21-09-2021"; B = $null;  C = -0.63; D = 1.15; E = 1.74; F = $null;  G = $null },
    @{ Row = 182; Date = "22-09-2021"; B = $null;  C = -0.61; D = 1.25; E = 1.77; F = 2.4;   G = 2.5 },
    @{ Row = 183; Date = "23-09-2021"; B = $null;  C = -0.44; D = 1.33; E = 1.83; F = 2.39;  G = 2.5 },
    @{ Row = 184; Date = "24-09-2021"; B = $null;  C = -0.39; D = 1.44; E = 1.97; F = 2.38;  G = $null },
    @{ Row = 185; Date = "27-09-2021"; B = $null;  C = $null;  D = 1.61; E = 2.12; F = 2.49;  G = $null },
    @{ Row = 186; Date = "28-09-2021"; B = -0.82; C = -0.45; D = 1.61; E = 2.14; F = 2.52;  G = 2.63 },
    @{ Row = 187; Date = "29-09-2021"; B = $null;  C = 0.09;  D = 1.66; E = 2.16; F = 2.58;  G = 2.66 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    if ($null -ne $r.B) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($null -ne $r.C) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($null -ne $r.D) { $ws.Cells.Item($r.Row, 4).Value = $r.D }
    if ($null -ne $r.E) { $ws.Cells.Item($r.Row, 5).Value = $r.E }
    if ($null -ne $r.F) { $ws.Cells.Item($r.Row, 6).Value = $r.F }
    if ($null -ne $r.G) { $ws.Cells.Item($r.Row, 7).Value = $r.G }
}
